# Update the "average" column for doctor_MA (column AF) for rows 4-13
# with the new results as described in the commit "updated results and code".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 0.65
    5  = 0.8
    6  = 0.717
    7  = 0.765
    8  = 0.8
    9  = 0.8
    10 = 0.8
    11 = 0.8
    12 = 1
    13 = 1.8
}

foreach ($row in $updates.Keys) {
    $ws.Range("AF$row").Value = $updates[$row]
}
